$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -2.74
$ws.Range("C3").Value = 4.69
$ws.Range("C4").Value = 4.69
$ws.Range("C5").Value = -0.28
$ws.Range("C6").Value = 0.11
$ws.Range("C7").Value = 5.73
$ws.Range("C8").Value = 1.17
$ws.Range("C9").Value = 0.38
$ws.Range("C11").Value = 2.14
$ws.Range("C12").Value = 1.38
$ws.Range("C13").Value = 1.1
$ws.Range("C14").Value = 2.46
$ws.Range("C15").Value = 1.9
$ws.Range("C16").Value = 0.49
$ws.Range("C17").Value = 2.46
$ws.Range("C18").Value = 2.46
$ws.Range("C19").Value = 2.22
$ws.Range("C20").Value = 2.17
$ws.Range("C21").Value = 1.57
$ws.Range("C22").Value = 2.49
$ws.Range("C23").Value = 2.54
$ws.Range("C24").Value = 2.11
$ws.Range("C25").Value = 2.76
$ws.Range("C26").Value = 3.09
$ws.Range("C27").Value = 2.29
$ws.Range("C28").Value = 3.18
$ws.Range("C29").Value = 3.1
$ws.Range("C30").Value = 2.92
$ws.Range("C31").Value = 3.07
$ws.Range("C32").Value = 3.14
$ws.Range("C33").Value = 3.06
$ws.Range("C34").Value = 3.08
$ws.Range("C36").Value = 2.95
$ws.Range("C37").Value = 3.26
$ws.Range("C38").Value = 2.38
$ws.Range("C39").Value = 3.03
$ws.Range("C40").Value = 3.47
$ws.Range("C41").Value = 3.63
$ws.Range("C42").Value = 2.25
$ws.Range("C43").Value = 3.67
$ws.Range("C44").Value = 2.44
$ws.Range("C45").Value = 3.17
$ws.Range("C46").Value = 3.84
$ws.Range("C47").Value = 3.93
$ws.Range("C48").Value = 4.2
$ws.Range("C49").Value = 3.17
$ws.Range("C50").Value = 3.32
$ws.Range("C51").Value = 3.64
$ws.Range("C52").Value = 3.38
$ws.Range("C53").Value = 3.85
$ws.Range("C54").Value = 3.42
$ws.Range("C55").Value = 3.93
$ws.Range("C56").Value = 2.25
$ws.Range("C57").Value = 4.33
$ws.Range("C58").Value = 4.47
$ws.Range("C59").Value = 3.5
$ws.Range("C60").Value = 3.58
$ws.Range("C61").Value = 5.09
$ws.Range("C62").Value = 4.68
$ws.Range("C63").Value = 3.59
$ws.Range("C64").Value = 3.96
$ws.Range("C65").Value = 4.23
$ws.Range("C66").Value = 5.29
$ws.Range("C67").Value = 4.33
$ws.Range("C68").Value = 4.44
$ws.Range("C69").Value = 4.06
$ws.Range("C70").Value = 3.31
$ws.Range("C71").Value = 4.69
$ws.Range("C72").Value = 4.07
$ws.Range("C73").Value = 4.89
$ws.Range("C74").Value = 5.05
$ws.Range("C75").Value = 4.76
$ws.Range("C76").Value = 5.05
$ws.Range("C77").Value = 4.44
$ws.Range("C78").Value = 4.8
$ws.Range("C79").Value = 4.16
$ws.Range("C80").Value = 4.62
$ws.Range("C81").Value = 5.09
$ws.Range("C82").Value = 4.98
$ws.Range("C83").Value = 5.61
$ws.Range("C84").Value = 4.36
$ws.Range("C85").Value = 5.34
$ws.Range("C86").Value = 4.96
$ws.Range("C88").Value = 5.65
$ws.Range("C89").Value = 5.34
$ws.Range("C91").Value = 3.74
$ws.Range("C92").Value = 5.73
$ws.Range("C93").Value = 5.52
$ws.Range("C94").Value = 5.97
$ws.Range("C96").Value = 6.34
$ws.Range("C97").Value = 6.07
$ws.Range("C98").Value = 6.29
$ws.Range("C100").Value = 8.54
$ws.Range("C101").Value = 8.04
